$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.05120575428009
$ws.Range("B1").Value = 1.176511883735657
$ws.Range("D1").Value = 1.650520205497742
$ws.Range("E1").Value = 0.9961546063423157
